# Demoting some runs from "current" status, and refreshing the sheet's
# filter range / view to include the newly-added J/K columns.
# https://app.asana.com/0/1201809392759895/1204436956308953/f

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# --- Clear the "status" (column G) value for the runs being demoted from
#     "current" back to blank. ---
$demotedRows = 18,31,47,54,55,56,57,63,64,65,66,77,83,88
foreach ($r in $demotedRows) {
    $ws.Range("G" + $r).Value = ""
}

# --- The used range now spans through column K, so refresh the AutoFilter
#     (and the _FilterDatabase defined name that backs it) to A1:K113. ---
$ws.Range("A1:K113").AutoFilter() | Out-Null

$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=all_runs!`$A`$1:`$K`$113"

# --- Narrow columns C and D now that fewer/shorter labels are shown. ---
$ws.Columns.Item(3).ColumnWidth = 42.14
$ws.Columns.Item(4).ColumnWidth = 12.29

# --- Update the frozen-pane scroll position and active selection. ---
$ws.Range("F10").Select() | Out-Null
